$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume(1h) (E) columns keep their text
# representation (values such as "1.011" or "27.000.15" must not be
# reinterpreted by Excel as numbers/dates).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.000.15"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.848.67"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "1.010"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "309.14"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("D8").Value = "0.3677"
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("D9").Value = "0.07226"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").Value = "0.9308"
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "0.07734"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "1.888.22"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "5.346"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "6.438"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "89.09"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "1.013"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "0.000008635"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "1.011"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "27.033.01"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("D22").Value = "5.070"
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "1.930"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").Value = "152.73"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "2.013"
$ws.Range("D28").Value = "114.20"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "4.962"
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("D30").Value = "0.08852"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "3.307"
$ws.Range("E31").Value = "  +4.89%  "
$ws.Range("D32").Value = "1.180"
$ws.Range("E32").Value = "  +0.29%  "
$ws.Range("D33").Value = "0.7399"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "4.500"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").Value = "2.728"
$ws.Range("E35").Value = "  -3.84%  "
$ws.Range("D36").Value = "1.109"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").Value = "0.01961"
$ws.Range("E37").Value = "  +2.11%  "
$ws.Range("D38").Value = "0.05263"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").Value = "2.974"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "0.5253"
$ws.Range("E40").Value = "  +3.74%  "
$ws.Range("D41").Value = "7.016"
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "0.1516"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").Value = "8.247"
$ws.Range("E43").Value = "  +1.69%  "
$ws.Range("D44").Value = "10.63"
$ws.Range("E44").Value = "  +4.69%  "
$ws.Range("D45").Value = "0.4749"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").Value = "1.011"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "101.77"
$ws.Range("D48").Value = "1.605"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "65.74"
$ws.Range("E49").Value = "  +2.59%  "
$ws.Range("D50").Value = "0.06066"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "0.8891"
$ws.Range("E51").Value = "  +3.97%  "
